$wb = $excel.ActiveWorkbook

# 1. Delete the "tblMode" worksheet entirely (table + all its rows go away).
$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("tblMode").Delete()
$excel.DisplayAlerts = $true

# 2. Rename a few worksheets to prefix them with "B " (they belong to the
#    "B" group of tables alongside "B tblJeu" / "B tblThemeJeu").
$wb.Worksheets.Item("tblVersion").Name = "B tblVersion"
$wb.Worksheets.Item("tblJeuSemblable").Name = "B tblJeuSemblable"
$wb.Worksheets.Item("tblPlateformeJeu").Name = "B tblPlateformeJeu"

# 3. Update the selection on "C tblSysExp" to a single cell (C85).
$ws = $wb.Worksheets.Item("C tblSysExp")
$ws.Activate()
[void]$ws.Range("C85").Select()

# 4. Update the selection on "B tblJeuSemblable" (was tblJeuSemblable).
$ws = $wb.Worksheets.Item("B tblJeuSemblable")
$ws.Activate()
[void]$ws.Range("K32").Select()

# 5. Update the selection on "B tblPlateformeJeu" and make it the active
#    sheet/tab (it becomes the last, focused tab in the saved workbook).
$ws = $wb.Worksheets.Item("B tblPlateformeJeu")
$ws.Activate()
[void]$ws.Range("E3").Select()
